# Fruta / hortaliza, semanal
# Insert a new weekly record row after the existing row 252 (pushing rows
# 253:338 down to 254:339) and populate the two affected rows with their
# correct data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 253; this shifts rows 253-338 down to
# 254-339 (data, styles, everything moves down), growing the sheet from
# 338 to 339 data rows.
$ws.Rows.Item(253).Insert()

# The row that used to be row 252 (date 2021-09-03, Primera, 250/2200/2600/2400, 800)
# now logically belongs at row 253 (it's still the same record, just renumbered).
# The newly inserted row 253 is currently empty, so fill it with that record.
$ws.Cells.Item(253, 1).Value2  = 9
$ws.Cells.Item(253, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(253, 3).Value2  = "Metropolitana"
$ws.Cells.Item(253, 4).Value2  = 44442
$ws.Cells.Item(253, 5).Value2  = 13
$ws.Cells.Item(253, 6).Value2  = 100112039
$ws.Cells.Item(253, 7).Value2  = "Ciboulette"
$ws.Cells.Item(253, 8).Value2  = "Sin especificar"
$ws.Cells.Item(253, 9).Value2  = "Primera"
$ws.Cells.Item(253, 10).Value2 = 250
$ws.Cells.Item(253, 11).Value2 = 2200
$ws.Cells.Item(253, 12).Value2 = 2600
$ws.Cells.Item(253, 13).Value2 = 2400
$ws.Cells.Item(253, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(253, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(253, 16).Value2 = 800
$ws.Cells.Item(253, 17).Value2 = 3
$ws.Cells.Item(253, 18).Value2 = "Hortaliza"
$ws.Cells.Item(253, 4).NumberFormat = $ws.Cells.Item(252, 4).NumberFormat

# Row 252 now carries a brand-new weekly record (date 2022-01-27) that was
# added by this edit; overwrite the price/volume fields accordingly
# (Mercado/Región/Categoría/Variedad/Calidad/Unidad/Origen stay the same).
$ws.Cells.Item(252, 4).Value2  = 44588
$ws.Cells.Item(252, 10).Value2 = 160
$ws.Cells.Item(252, 11).Value2 = 800
$ws.Cells.Item(252, 12).Value2 = 1000
$ws.Cells.Item(252, 13).Value2 = 900
$ws.Cells.Item(252, 16).Value2 = 300
